# Project Optimization.xlsx - "Changed tcl file to correct freq"
#
# The underlying edit appends a new measured data point (row 33) to the
# bottom of the second data table (rows 22-31), leaving a blank spacer row
# (32) ahead of it -- mirroring the existing spacer row (21) that separates
# table 1 (rows 2-20) from table 2 (rows 22-31). It also clears a stray
# leftover formula in that older spacer row (T21) that was evaluating to
# #DIV/0!.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the stray #DIV/0! formula left over in the separator row (21) ---
# T21 used to hold a copy of the shared formula P21*I21/D21 which divides by
# the (blank -> 0) D21, producing #DIV/0!. The row is otherwise empty, so the
# formula is simply removed, leaving the cell blank (format untouched).
$ws.Range("T21").ClearContents()

# --- 2. New row 32: blank spacer row, but the fill-down formulas for
#        columns C/D/N/O/P/T continue into it (C/D evaluate to 0 because
#        B32 is blank; N/O/P/T are left as blank-but-formatted cells, same
#        as the equivalent cells on spacer row 21) ---
$ws.Range("C32").Formula = "=B32*10^3"
$ws.Range("D32").Formula = "=B32*10^9"
$ws.Range("N21:P21").Copy()
$ws.Range("N32:P32").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("T21").Copy()
$ws.Range("T32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. New row 33: a new measured data point appended to table 2 ---
$ws.Range("A33").Value = 0.28999999999999998
$ws.Range("B33").Value = 3.448275862
$ws.Range("C33").Formula = "=B33*10^3"
$ws.Range("D33").Formula = "=B33*10^9"
$ws.Range("E33").Value = 1588.654166
$ws.Range("F33").Value = 71.668609000000004
$ws.Range("G33").Value = 1919.2955489999999
$ws.Range("H33").Value = 759.86486100000002
$ws.Range("I33").Value = 4267.8145750000003
$ws.Range("J33").Value = 87.594999999999999
$ws.Range("K33").Value = 7900
$ws.Range("L33").Value = 7930000000
$ws.Range("M33").Value = 15900
$ws.Range("N33").Formula = "=J33+K33"
$ws.Range("O33").Formula = "=N33*10^6"
$ws.Range("P33").Formula = "=O33/D33"
$ws.Range("Q33").Value = 0.06
$ws.Range("T33").Formula = "=P33*I33/D33"

# --- 4. Formatting: mirror the formats used by the rest of table 2 (rows
#        22-31) onto the two new rows so the appended data reads the same
#        as the existing rows (number formats, borders already established
#        by the surrounding cells). Scoped to the exact columns that carry
#        data so no stray empty cells are introduced in unused columns. ---
$ws.Range("A31:Q31").Copy()
$ws.Range("A33:Q33").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("T31").Copy()
$ws.Range("T33").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C31:D31").Copy()
$ws.Range("C32:D32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 5. Move the selection to match where the author was last working ---
$ws.Range("T33").Select()
